# Update MSME Slovenia summary figures with refined (more precise) values.
# Each target cell currently holds its number formatted as plain text
# (shared string), so we force the new value to stay text too -- otherwise
# a numeric-looking string gets auto-converted to a real number by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "B11" "36.29"
Set-TextValue "C11" "30.67"
Set-TextValue "D11" "66.96"

Set-TextValue "D12" "24.71"

Set-TextValue "B33" "53.75"
Set-TextValue "C33" "3.37"
Set-TextValue "D33" "57.12"

Set-TextValue "B34" "32.93"
Set-TextValue "D34" "72.92"

Set-TextValue "B36" "93.92"
Set-TextValue "C36" "5.88"
Set-TextValue "D36" "99.81"

Set-TextValue "C40" "42.36"
Set-TextValue "D40" "63.56"
